# Added New Mac-Address and Document Types
# Appends 5 new device rows (Finger Print Scanner 32, IRIS Scanner 32,
# Web Camera 32, Document Scanner 32, Printer 32) to the device master
# table on Sheet1, plus 5 trailing formatted-but-empty rows, matching the
# existing row layout (columns A-K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlHAlignLeft - reproduces the existing "style index 1" (left aligned)
# cell format used throughout column H (the is_active boolean column).
$xlHAlignLeft = -4131

# New data rows to append, in order: row, id, name, mac_address,
# serial_num, dspec_id
$newRows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; Dspec = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";          Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; Dspec = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";            Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; Dspec = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";      Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; Dspec = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";                Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; Dspec = 920 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Id
    $ws.Range("B$row").Value = $r.Name
    $ws.Range("C$row").Value = $r.Mac
    $ws.Range("D$row").Value = $r.Serial
    $ws.Range("F$row").Value = $r.Dspec
    $ws.Range("G$row").Value = "eng"
    $ws.Range("H$row").Value = $true
    $ws.Range("H$row").HorizontalAlignment = $xlHAlignLeft
    $ws.Range("I$row").Value = "superadmin"
    $ws.Range("J$row").Value = "now()"
    $ws.Range("K$row").Value = "now()"
}

# Five trailing rows (162-166) that only carry the left-aligned format on
# column H, with no values - mirrors the source workbook exactly.
for ($row = 162; $row -le 166; $row++) {
    $ws.Range("H$row").HorizontalAlignment = $xlHAlignLeft
}

# Update the visible window/selection to match where the author ended up
# after typing in the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 154
$win.ScrollColumn = 1
$ws.Range("E159").Select()

Write-Host "Added 5 device rows (157-161) and 5 formatted rows (162-166)"
